$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (border/bold/center/top alignment) from E1 to F1, then set header text
$ws.Range("E1").Copy($ws.Range("F1"))
$ws.Range("F1").Value = "time_taken"

# Fill in the time_taken values for each data row (F2:F66)
$timeTaken = @(
    "2021-10-05 10:50:46.694758",
    "2021-10-05 10:50:46.694770",
    "2021-10-05 10:50:46.694773",
    "2021-10-05 10:50:46.694776",
    "2021-10-05 10:50:46.694779",
    "2021-10-05 10:50:46.694781",
    "2021-10-05 10:50:46.694784",
    "2021-10-05 10:50:46.694786",
    "2021-10-05 10:50:46.694789",
    "2021-10-05 10:50:46.694792",
    "2021-10-05 10:50:46.694794",
    "2021-10-05 10:50:46.694797",
    "2021-10-05 10:50:46.694799",
    "2021-10-05 10:50:46.694802",
    "2021-10-05 10:50:46.694804",
    "2021-10-05 10:50:46.694807",
    "2021-10-05 10:50:46.694809",
    "2021-10-05 10:50:46.694812",
    "2021-10-05 10:50:46.694815",
    "2021-10-05 10:50:46.694817",
    "2021-10-05 10:50:46.694820",
    "2021-10-05 10:50:46.694822",
    "2021-10-05 10:50:46.694824",
    "2021-10-05 10:50:46.694827",
    "2021-10-05 10:50:46.694830",
    "2021-10-05 10:50:46.694832",
    "2021-10-05 10:50:46.694835",
    "2021-10-05 10:50:46.694837",
    "2021-10-05 10:50:46.694840",
    "2021-10-05 10:50:46.694842",
    "2021-10-05 10:50:46.694845",
    "2021-10-05 10:50:46.694847",
    "2021-10-05 10:50:46.694850",
    "2021-10-05 10:50:46.694853",
    "2021-10-05 10:50:46.694855",
    "2021-10-05 10:50:46.694858",
    "2021-10-05 10:50:46.694860",
    "2021-10-05 10:50:46.694863",
    "2021-10-05 10:50:46.694865",
    "2021-10-05 10:50:46.694868",
    "2021-10-05 10:50:46.694871",
    "2021-10-05 10:50:46.694873",
    "2021-10-05 10:50:46.694876",
    "2021-10-05 10:50:46.694878",
    "2021-10-05 10:50:46.694881",
    "2021-10-05 10:50:46.694883",
    "2021-10-05 10:50:46.694886",
    "2021-10-05 10:50:46.694888",
    "2021-10-05 10:50:46.694891",
    "2021-10-05 10:50:46.694893",
    "2021-10-05 10:50:46.694896",
    "2021-10-05 10:50:46.694898",
    "2021-10-05 10:50:46.694901",
    "2021-10-05 10:50:46.694904",
    "2021-10-05 10:50:46.694906",
    "2021-10-05 10:50:46.694909",
    "2021-10-05 10:50:46.694911",
    "2021-10-05 10:50:46.694914",
    "2021-10-05 10:50:46.694916",
    "2021-10-05 10:50:46.694919",
    "2021-10-05 10:50:46.694921",
    "2021-10-05 10:50:46.694924",
    "2021-10-05 10:50:46.694926",
    "2021-10-05 10:50:46.694929",
    "2021-10-05 10:50:46.694932"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timeTaken[$i]
}

Write-Output "Applied time_taken column to F1:F66"
